# Auto-generated edit script applying market-data value updates
# as captured by the authoritative diff (scheduled runner sheet sync).
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2588.1765
$ws.Range("I40").Value = 1399.9166
$ws.Range("J40").Value = 5440
$ws.Range("K40").Value = 1399.9166
$ws.Range("L40").Value = 5440
$ws.Range("M40").Value = -1224.9166
$ws.Range("N40").Value = -5790
$ws.Range("H113").Value = 8160.923
$ws.Range("J113").Value = 9141.857
$ws.Range("L113").Value = 9141.857
$ws.Range("N113").Value = -15649.857
$ws.Range("H137").Value = 4897.4546
$ws.Range("I137").Value = 4718.25
$ws.Range("K137").Value = 14154.75
$ws.Range("M137").Value = -11604.75

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H10").Value = 0
$ws.Range("I10").Value = 0
$ws.Range("K10").Value = 0
$ws.Range("M10").ClearContents()
$ws.Range("H46").Value = 3576
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 3576
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 3576
$ws.Range("M46").ClearContents()
$ws.Range("N46").Value = -4214
$ws.Range("H97").Value = 657.1429000000001
$ws.Range("I97").Value = 657.1429000000001
$ws.Range("K97").Value = 657.1429000000001
$ws.Range("M97").Value = -161.1429000000001

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("M15").ClearContents()
$ws.Range("N15").ClearContents()
$ws.Range("H94").Value = 2452.7
$ws.Range("I94").Value = 2452.7
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 2452.7
$ws.Range("L94").Value = 0
$ws.Range("M94").Value = -2001.7
$ws.Range("N94").ClearContents()
$ws.Range("H105").Value = 1999.8
$ws.Range("I105").Value = 1999
$ws.Range("J105").Value = 1999.9231
$ws.Range("K105").Value = 1999
$ws.Range("L105").Value = 1999.9231
$ws.Range("M105").Value = -252
$ws.Range("N105").Value = -5493.9231
$ws.Range("H107").Value = 712.2857
$ws.Range("I107").Value = 619.1111
$ws.Range("K107").Value = 619.1111
$ws.Range("M107").Value = 1300.8889

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1476
$ws.Range("I5").Value = 1073.25
$ws.Range("J5").Value = 1798.2
$ws.Range("K5").Value = 3219.75
$ws.Range("L5").Value = 5394.6
$ws.Range("M5").Value = -3107.75
$ws.Range("N5").Value = -5618.6
$ws.Range("H32").Value = 893.6667
$ws.Range("J32").Value = 893.6667
$ws.Range("L32").Value = 2681.0001
$ws.Range("N32").Value = -3247.0001
$ws.Range("H35").Value = 1099
$ws.Range("I35").Value = 1099
$ws.Range("K35").Value = 3297
$ws.Range("M35").Value = -3009
$ws.Range("H37").Value = 99997.75
$ws.Range("J37").Value = 99997.75
$ws.Range("L37").Value = 299993.25
$ws.Range("N37").Value = -300217.25
$ws.Range("H68").Value = 999
$ws.Range("J68").Value = 998.75
$ws.Range("L68").Value = 2996.25
$ws.Range("N68").Value = -4618.25
$ws.Range("H71").Value = 999
$ws.Range("J71").Value = 998.75
$ws.Range("L71").Value = 8988.75
$ws.Range("N71").Value = -17100.75
$ws.Range("H135").Value = 1476
$ws.Range("I135").Value = 1073.25
$ws.Range("J135").Value = 1798.2
$ws.Range("K135").Value = 9659.25
$ws.Range("L135").Value = 16183.8
$ws.Range("M135").Value = -7124.25
$ws.Range("N135").Value = -21253.8

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 47625252
$ws.Range("I70").Value = 47625252
$ws.Range("K70").Value = 47625252
$ws.Range("M70").Value = -47624982
$ws.Range("H73").Value = 47625252
$ws.Range("I73").Value = 47625252
$ws.Range("K73").Value = 47625252
$ws.Range("M73").Value = -47624316
$ws.Range("H102").Value = 3196.8
$ws.Range("I102").Value = 3188.3076
$ws.Range("J102").Value = 3252
$ws.Range("K102").Value = 3188.3076
$ws.Range("L102").Value = 3252
$ws.Range("M102").Value = -1566.3076
$ws.Range("N102").Value = -6496
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2321.8076
$ws.Range("I46").Value = 1948.0625
$ws.Range("K46").Value = 1948.0625
$ws.Range("M46").Value = -1760.0625
$ws.Range("H61").Value = 4279.6
$ws.Range("I61").Value = 3799.6667
$ws.Range("J61").Value = 4999.5
$ws.Range("K61").Value = 3799.6667
$ws.Range("L61").Value = 4999.5
$ws.Range("M61").Value = -3597.6667
$ws.Range("N61").Value = -5403.5
$ws.Range("H93").Value = 4499
$ws.Range("I93").Value = 4499
$ws.Range("K93").Value = 4499
$ws.Range("M93").Value = -3251
$ws.Range("H113").Value = 4279.6
$ws.Range("I113").Value = 3799.6667
$ws.Range("J113").Value = 4999.5
$ws.Range("K113").Value = 3799.6667
$ws.Range("L113").Value = 4999.5
$ws.Range("M113").Value = -1629.6667
$ws.Range("N113").Value = -9339.5
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("N122").ClearContents()

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 82162.5
$ws.Range("I5").Value = 88200
$ws.Range("J5").Value = 51975
$ws.Range("K5").Value = 88200
$ws.Range("L5").Value = 88200
$ws.Range("M5").Value = -88088
$ws.Range("N5").Value = -52199
$ws.Range("H70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").ClearContents()
$ws.Range("H73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").ClearContents()
$ws.Range("H96").Value = 4000
$ws.Range("J96").Value = 4000
$ws.Range("L96").Value = 4000
$ws.Range("N96").Value = -6746
$ws.Range("H107").Value = 1617
$ws.Range("I107").Value = 1028.3334
$ws.Range("K107").Value = 3085.0002
$ws.Range("M107").Value = -1165.0002
$ws.Range("H122").Value = 1500
$ws.Range("I122").Value = 1500
$ws.Range("K122").Value = 4500
$ws.Range("M122").Value = -2050
$ws.Range("H131").Value = 0
$ws.Range("J131").Value = 0
$ws.Range("L131").Value = 0
$ws.Range("N131").ClearContents()
